$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The edit adds three new command rows to the table:
#   - "Set RSSI Monitor interval"  -> becomes new row 28
#   - "Get LiPo Voltage"           -> becomes new row 30
#   - "Get single RSSI reading"    -> becomes new row 31
# (the former row 28 "Skip/Enable First Lap" keeps its content and becomes
#  row 29; everything from the old row 29 onward shifts down accordingly)
# and appends one blank trailing row at the end (row 37).
# ---------------------------------------------------------------------------

# Insert 3 blank rows at the target positions (tracking the downward shift):
#  - row 28 is brand new
#  - row 29 keeps the old row 28 content ("Skip/Enable First Lap")
#  - rows 30 and 31 are brand new
#  - old row 29 onward shifts down to 32+
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(30).Insert()
$ws.Rows.Item(30).Insert()

# ---------------------------------------------------------------------------
# Row 28: "Set RSSI Monitor interval"
# ---------------------------------------------------------------------------
$ws.Range("A28").Value = "Set RSSI Monitor interval"
$ws.Range("B28").Value = 2
$ws.Range("C28").Value = "d"
$ws.Range("D28").Value = "d"
$ws.Range("E28").Value = "RSSI monitor interval to set (2 bytes). Values are measured in ~0.1ms intervals. e.g. value of 1000 (decimal) gives 100ms. Min allowed value is 0xA."
$ws.Range("F28").Value = "RSSI monitor interval (2 bytes) (see request description)"
$ws.Range("G28").Value = "R1d0001\n -> S1d000A\n"

# ---------------------------------------------------------------------------
# Row 30: "Get LiPo Voltage"
# ---------------------------------------------------------------------------
$ws.Range("A30").Value = "Get LiPo Voltage"
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = "Y"
$ws.Range("D30").Value = "Y"
$ws.Range("F30").Value = "Reading from analog voltage pin (2 bytes). Response comes only from nodes that have LiPo monitoring."
$ws.Range("G30").Value = "R*Y\n -> S1Y00D5\n"

# ---------------------------------------------------------------------------
# Row 31: "Get single RSSI reading"
# ---------------------------------------------------------------------------
$ws.Range("A31").Value = "Get single RSSI reading"
$ws.Range("B31").Value = 2
$ws.Range("C31").Value = "E"
$ws.Range("D31").Value = "S"
$ws.Range("F31").Value = "Current RSSI Value (2 bytes)"
$ws.Range("G31").Value = "R0E\n -> S0S0111\n"

# ---------------------------------------------------------------------------
# Re-apply the table's standard formatting (borders/fonts/fills/alignment)
# to the new rows, by copying the format from an existing row that already
# uses the same column style pattern (row 14: A/E/F left-aligned Consolas,
# B quoted Segoe UI, C/D centered Consolas, G left-aligned Consolas).
# This is done after the values are set so that numeric cells in column B
# keep the "quote-prefix" cell style used throughout the rest of the table.
# ---------------------------------------------------------------------------
$ws.Rows.Item(14).Copy()
$ws.Range("A28:G28").PasteSpecial(-4122)
$ws.Rows.Item(14).Copy()
$ws.Range("A30:G30").PasteSpecial(-4122)
$ws.Rows.Item(14).Copy()
$ws.Range("A31:G31").PasteSpecial(-4122)

# Rows 30 and 31 have no content/description in column E, so remove the
# empty cell entirely (it should not be present at all in the final sheet).
$ws.Range("E30").Clear()
$ws.Range("E31").Clear()

# Row heights for the new rows
$ws.Rows.Item(28).RowHeight = 75
$ws.Rows.Item(30).RowHeight = 48.75
$ws.Rows.Item(31).RowHeight = 48.75

# ---------------------------------------------------------------------------
# Append the trailing blank row 37
# ---------------------------------------------------------------------------
$ws.Rows.Item(37).RowHeight = 31.5

# ---------------------------------------------------------------------------
# Update the selection to match the author's final cursor position
# ---------------------------------------------------------------------------
$ws.Rows.Item(32).Select()
